$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.159713
$ws.Range("H2").Value = 0.319426
$ws.Range("M2").Value = 0.8108664999999999
$ws.Range("N2").Value = 1.621733
$ws.Range("O2").Value = 0.0929863610355799
$ws.Range("P2").Value = 0.08514801891941343
$ws.Range("Q2").Value = 0.1295059213145
$ws.Range("R2").Value = 0.5180236852579999
$ws.Range("S2").Value = 0.0929863610355799
$ws.Range("T2").Value = 0.08514801891941343

# Row 3
$ws.Range("G3").Value = 0.159713
$ws.Range("H3").Value = 0.319426
$ws.Range("O3").Value = 0.1736650041232461
$ws.Range("P3").Value = 0.2385387097426263
$ws.Range("Q3").Value = 0.2418703787156667
$ws.Range("R3").Value = 1.451222272294
$ws.Range("S3").Value = 0.1736650041232461
$ws.Range("T3").Value = 0.2385387097426263

# Row 4
$ws.Range("G4").Value = 0.159713
$ws.Range("H4").Value = 0.319426
$ws.Range("M4").Value = 0.01802466666666666
$ws.Range("N4").Value = 0.054074
$ws.Range("O4").Value = 0.002066984099371042
$ws.Range("P4").Value = 0.002839119617747411
$ws.Range("Q4").Value = 0.002878773587333333
$ws.Range("R4").Value = 0.017272641524
$ws.Range("S4").Value = 0.002066984099371042
$ws.Range("T4").Value = 0.002839119617747411

# Row 5
$ws.Range("G5").Value = 0.159713
$ws.Range("H5").Value = 0.319426
$ws.Range("M5").Value = 6.303909
$ws.Range("N5").Value = 12.607818
$ws.Range("O5").Value = 0.7229026704265641
$ws.Range("P5").Value = 0.6619651481449297
$ws.Range("Q5").Value = 1.006816218117
$ws.Range("R5").Value = 4.027264872468
$ws.Range("S5").Value = 0.7229026704265641
$ws.Range("T5").Value = 0.6619651481449297

# Row 6
$ws.Range("G6").Value = 0.159713
$ws.Range("H6").Value = 0.319426
$ws.Range("M6").Value = 0.07084766666666666
$ws.Range("N6").Value = 0.212543
$ws.Range("O6").Value = 0.008124477594271174
$ws.Range("P6").Value = 0.01115942968737079
$ws.Range("Q6").Value = 0.01131529338633333
$ws.Range("R6").Value = 0.067891760318
$ws.Range("S6").Value = 0.008124477594271174
$ws.Range("T6").Value = 0.01115942968737079

# Row 7
$ws.Range("G7").Value = 0.159713
$ws.Range("H7").Value = 0.319426
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.002219333333333333
$ws.Range("N7").Value = 0.006658
$ws.Range("O7").Value = 0.0002545027209677923
$ws.Range("P7").Value = 0.0003495738879121623
$ws.Range("Q7").Value = 0.0003544563846666666
$ws.Range("R7").Value = 0.002126738308
$ws.Range("S7").Value = 0.0002545027209677923
$ws.Range("T7").Value = 0.0003495738879121623
